$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55: politeness_score (B55) was stored as text "3"; fix it to a real number 3
$ws.Range("B55").Value = 3

# Row 56: new annotation row for Ying Tang
$ws.Range("A56").Value = "Ying Tang"

# B56 keeps the politeness_score as text "3" (matches source data pattern)
$ws.Range("B56").NumberFormat = "@"
$ws.Range("B56").Value = "3"
$ws.Range("B56").NumberFormat = "General"
$ws.Range("B56").ClearFormats()

$ws.Range("C56").Value = "This paper shows that"
$ws.Range("D56").Value = "SMY"
$ws.Range("E56").Value = "RES"
$ws.Range("F56").Value = "e3eeb88f-0832-4aa9-a6cc-39ada0451b32"
$ws.Range("G56").Value = "BkJ3ibb0-_annotated.xlsx"
$ws.Range("H56").Value = "This paper shows that models trained on a synthetic dataset are vulnerable to small adversarial perturbations which lie on the data manifold."
